$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "20.241.30"
Set-TextValue $ws.Range("E2") "  +1.90%  "
Set-TextValue $ws.Range("D3") "1.443.01"
Set-TextValue $ws.Range("E3") "  +3.88%  "
Set-TextValue $ws.Range("D4") "1.006"
Set-TextValue $ws.Range("E4") "  +0.15%  "
Set-TextValue $ws.Range("D5") "0.9138"
Set-TextValue $ws.Range("E5") "  -8.96%  "
Set-TextValue $ws.Range("D6") "278.21"
Set-TextValue $ws.Range("E6") "  +3.93%  "
Set-TextValue $ws.Range("D7") "0.3649"
Set-TextValue $ws.Range("E7") "  +1.01%  "
Set-TextValue $ws.Range("D8") "0.3113"
Set-TextValue $ws.Range("E8") "  +3.66%  "
Set-TextValue $ws.Range("D9") "39.16"
Set-TextValue $ws.Range("E9") "  +0.57%  "
Set-TextValue $ws.Range("D10") "1.020"
Set-TextValue $ws.Range("E10") "  +7.15%  "
Set-TextValue $ws.Range("D11") "0.06524"
Set-TextValue $ws.Range("E11") "  +3.13%  "
Set-TextValue $ws.Range("D12") "1.000"
Set-TextValue $ws.Range("E12") "  -0.19%  "
Set-TextValue $ws.Range("D13") "5.397"
Set-TextValue $ws.Range("E13") "  +4.05%  "
Set-TextValue $ws.Range("D14") "17.53"
Set-TextValue $ws.Range("E14") "  +8.22%  "
Set-TextValue $ws.Range("D15") "6.053"
Set-TextValue $ws.Range("E15") "  +1.44%  "
Set-TextValue $ws.Range("D16") "0.00001016"
Set-TextValue $ws.Range("E16") "  +4.00%  "
Set-TextValue $ws.Range("D17") "1.444.16"
Set-TextValue $ws.Range("E17") "  +3.52%  "
Set-TextValue $ws.Range("D18") "0.9443"
Set-TextValue $ws.Range("E18") "  -5.87%  "
Set-TextValue $ws.Range("D19") "0.05637"
Set-TextValue $ws.Range("E19") "  +0.77%  "
Set-TextValue $ws.Range("D20") "68.16"
Set-TextValue $ws.Range("E20") "  -1.49%  "
Set-TextValue $ws.Range("D23") "10.82"
Set-TextValue $ws.Range("E23") "  +3.91%  "
Set-TextValue $ws.Range("D24") "2.260"
Set-TextValue $ws.Range("E24") "  +0.85%  "
Set-TextValue $ws.Range("D25") "20.246.73"
Set-TextValue $ws.Range("E25") "  +1.82%  "
Set-TextValue $ws.Range("D26") "2.163"
Set-TextValue $ws.Range("E26") "  +2.56%  "
Set-TextValue $ws.Range("D27") "137.61"
Set-TextValue $ws.Range("E27") "  +2.55%  "
Set-TextValue $ws.Range("D28") "16.93"
Set-TextValue $ws.Range("E28") "  +3.75%  "
Set-TextValue $ws.Range("D29") "1.597.11"
Set-TextValue $ws.Range("E29") "  +2.82%  "
Set-TextValue $ws.Range("D30") "110.02"
Set-TextValue $ws.Range("E30") "  +3.74%  "
Set-TextValue $ws.Range("D31") "3.826"
Set-TextValue $ws.Range("E31") "  -0.06%  "
Set-TextValue $ws.Range("D32") "0.8015"
Set-TextValue $ws.Range("E32") "  +3.57%  "
Set-TextValue $ws.Range("D33") "4.840"
Set-TextValue $ws.Range("E33") "  -5.61%  "
Set-TextValue $ws.Range("D34") "0.07690"
Set-TextValue $ws.Range("E34") "  +2.39%  "
Set-TextValue $ws.Range("D35") "0.05925"
Set-TextValue $ws.Range("E35") "  +7.82%  "
Set-TextValue $ws.Range("D36") "1.447"
Set-TextValue $ws.Range("E36") "  +11.54%  "
Set-TextValue $ws.Range("D37") "4.672"
Set-TextValue $ws.Range("E37") "  +1.58%  "
Set-TextValue $ws.Range("D38") "1.133"
Set-TextValue $ws.Range("E38") "  +10.54%  "
Set-TextValue $ws.Range("D39") "0.01994"
Set-TextValue $ws.Range("E39") "  +1.07%  "
Set-TextValue $ws.Range("D40") "10.16"
Set-TextValue $ws.Range("E40") "  +3.70%  "
Set-TextValue $ws.Range("D41") "0.9320"
Set-TextValue $ws.Range("E41") "  -7.04%  "
Set-TextValue $ws.Range("D42") "0.1836"
Set-TextValue $ws.Range("E42") "  -1.28%  "
Set-TextValue $ws.Range("D43") "7.143"
Set-TextValue $ws.Range("E43") "  -12.47%  "
Set-TextValue $ws.Range("D44") "3.522"
Set-TextValue $ws.Range("E44") "  +2.21%  "
Set-TextValue $ws.Range("D45") "0.5222"
Set-TextValue $ws.Range("E45") "  +2.41%  "
Set-TextValue $ws.Range("D46") "11.99"
Set-TextValue $ws.Range("E46") "  +3.13%  "
Set-TextValue $ws.Range("D47") "118.94"
Set-TextValue $ws.Range("E47") "  +11.30%  "
Set-TextValue $ws.Range("D48") "0.5131"
Set-TextValue $ws.Range("E48") "  +5.12%  "
Set-TextValue $ws.Range("D49") "1.755"
Set-TextValue $ws.Range("E49") "  +4.11%  "
Set-TextValue $ws.Range("D50") "0.06337"
Set-TextValue $ws.Range("E50") "  +4.45%  "
Set-TextValue $ws.Range("D51") "0.9919"
Set-TextValue $ws.Range("E51") "  -0.99%  "

# Row 21/22 swap: coin name, link, and new price/volume values
Set-TextValue $ws.Range("B21") "Avalanche"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D21") "14.40"
Set-TextValue $ws.Range("E21") "  +2.83%  "

Set-TextValue $ws.Range("B22") "Uniswap"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "5.387"
Set-TextValue $ws.Range("E22") "  -0.53%  "
